# edit.ps1
# Applies the changes described by the diff between south_des_moines.xlsx
# (before) and the updated workbook (after):
#   1. Update the "Exported On:" timestamp rich-text run in A2.
#   2. Rename two "do not use" style items (B24, B179).
#   3. Update many Quantity (column D) values, including the Total (D183).
#   4. Slightly widen column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Update the exported-on timestamp, keeping the "Exported On: " bold
#    label run intact and only changing the date/time run's text.
# ---------------------------------------------------------------------
$dateCell = $ws.Range("A2")
$label = "Exported On: "
$dateCell.Characters($label.Length + 1, 20).Text = "06/08/2025 12:31 PM"
$dateCell.Characters(1, $label.Length).Font.Bold = $true
$dateCell.Characters($label.Length + 1, 20).Font.Bold = $false

# ---------------------------------------------------------------------
# 2. Rename "do not use" placeholder items.
# ---------------------------------------------------------------------
$ws.Range("B24").Value = "ZZZZZZZZ DO NOT USE"
$ws.Range("B179").Value = "Do not use"

# ---------------------------------------------------------------------
# 3. Update Quantity (column D) values.
# ---------------------------------------------------------------------
$ws.Range("D9").Value = 14.0
$ws.Range("D10").Value = 17.0
$ws.Range("D11").Value = 4.0
$ws.Range("D12").Value = 15.0
$ws.Range("D13").Value = 12.0
$ws.Range("D14").Value = 3.0
$ws.Range("D16").Value = 4.0
$ws.Range("D28").Value = 1.0
$ws.Range("D30").Value = 17.0
$ws.Range("D31").Value = 6.0
$ws.Range("D36").Value = 14.0
$ws.Range("D38").Value = 13.0
$ws.Range("D39").Value = 5.0
$ws.Range("D40").Value = 3.0
$ws.Range("D42").Value = 6.0
$ws.Range("D43").Value = 12.0
$ws.Range("D49").Value = 3.0
$ws.Range("D52").Value = 1.0
$ws.Range("D57").Value = 9.0
$ws.Range("D61").Value = 1.0
$ws.Range("D63").Value = 14.0
$ws.Range("D64").Value = 3.0
$ws.Range("D67").Value = 5.0
$ws.Range("D68").Value = 4.0
$ws.Range("D70").Value = 15.0
$ws.Range("D71").Value = 5.0
$ws.Range("D89").Value = 7.0
$ws.Range("D100").Value = 8.0
$ws.Range("D101").Value = 9.0
$ws.Range("D106").Value = 10.0
$ws.Range("D113").Value = 9.0
$ws.Range("D115").Value = 8.0
$ws.Range("D116").Value = 2.0
$ws.Range("D117").Value = 5.0
$ws.Range("D118").Value = 9.0
$ws.Range("D119").Value = 6.0
$ws.Range("D121").Value = 2.0
$ws.Range("D122").Value = 2.0
$ws.Range("D123").Value = 2.0
$ws.Range("D124").Value = 1.0
$ws.Range("D127").Value = 3.0
$ws.Range("D129").Value = 4.0
$ws.Range("D131").Value = 8.0
$ws.Range("D134").Value = 11.0
$ws.Range("D135").Value = 7.0
$ws.Range("D136").Value = 7.0
$ws.Range("D137").Value = 4.0
$ws.Range("D141").Value = 2.0
$ws.Range("D142").Value = 1.0
$ws.Range("D144").Value = 2.0
$ws.Range("D145").Value = 0.0
$ws.Range("D146").Value = 0.0
$ws.Range("D148").Value = 1.0
$ws.Range("D149").Value = 4.0
$ws.Range("D150").Value = 1.0
$ws.Range("D151").Value = 3.0
$ws.Range("D152").Value = 4.0
$ws.Range("D154").Value = 3.0
$ws.Range("D156").Value = 2.0
$ws.Range("D164").Value = 4.0
$ws.Range("D169").Value = 4.0
$ws.Range("D172").Value = 5.0
$ws.Range("D173").Value = 5.0
$ws.Range("D175").Value = 1.0
$ws.Range("D177").Value = 7.0
$ws.Range("D178").Value = 0.0
$ws.Range("D180").Value = 8.0
$ws.Range("D181").Value = 8.0
$ws.Range("D183").Value = 949.0

# ---------------------------------------------------------------------
# 4. Widen column A slightly (bestFit width nudged from 20.17 to 20.75).
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20.04
